$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column D
$ws.Range("D1").Value = "ITI"

# Update column C (ConditionType) and add column D (ITI) values for rows 2-17
$data = @(
    @(3, 6),
    @(1, 6),
    @(2, 7),
    @(3, 7),
    @(4, 6),
    @(2, 8),
    @(4, 7),
    @(1, 6),
    @(4, 6),
    @(4, 6),
    @(1, 8),
    @(1, 6),
    @(3, 7),
    @(3, 7),
    @(2, 8),
    @(2, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
}

# Delete the old rows 18-20 (trials 17, 18, 19) which are no longer present
$ws.Range("A18:D20").Delete()

# Update selection to match target state
$ws.Range("C8").Select()
